# Applies updated First_Noticeable_Increase_Index (C), First_Noticeable_Increase_Cumulative_Value (E)
# and Pulse_Width (G) values to the Step3_DataPts_* sheets, reflecting the newly introduced
# zero_before_threshold parameter in the upstream signal-processing tool.

$wb = $excel.ActiveWorkbook

# Values are identical for column C and E across all four threshold sheets (0.5 / 0.7 / 0.8 / 0.9)
# because they are threshold-independent; only Pulse_Width (G) differs per sheet.
$commonRows = @(
    @{ Row = 2; C = 87; E = 0.006588967598699177 },
    @{ Row = 3; C = 87; E = 0.002549111923378791 },
    @{ Row = 4; C = 87; E = 0.003116462355016064 },
    @{ Row = 5; C = 87; E = 0.00264633726401236  },
    @{ Row = 6; C = 88; E = 0.00339759507480527  }
)

$sheetPulseWidths = @{
    "Step3_DataPts_0.5" = @{ 2 = 15; 3 = 6;  4 = 8;  5 = 6;  6 = 6  }
    "Step3_DataPts_0.7" = @{ 2 = 42; 3 = 31; 4 = 34; 5 = 31; 6 = 31 }
    "Step3_DataPts_0.8" = @{ 2 = 64; 3 = 51; 4 = 56; 5 = 51; 6 = 45 }
    "Step3_DataPts_0.9" = @{ 2 = 77; 3 = 74; 4 = 75; 5 = 74; 6 = 72 }
}

foreach ($sheetName in $sheetPulseWidths.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $pulseWidths = $sheetPulseWidths[$sheetName]

    foreach ($entry in $commonRows) {
        $row = $entry.Row
        $ws.Cells.Item($row, 3).Value = $entry.C
        $ws.Cells.Item($row, 5).Value = $entry.E
        $ws.Cells.Item($row, 7).Value = $pulseWidths[$row]
    }
}
